# Refresh the crypto price / 1h-volume-change snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): the column stores plain text (e.g. '34.516.96',
#     '0.0967') rather than numbers. A leading apostrophe forces Excel to
#     keep each new value as a literal string instead of auto-parsing it as a
#     number/date (which would also silently drop meaningful trailing zeros,
#     e.g. '0.600' -> 0.6). The apostrophe marks the cell with a 'quote
#     prefix' style, so afterwards we reset those cells back to the default
#     'Normal' style to match the original (unstyled) formatting.
$ws.Range("D2").Value = "'34.516.96"
$ws.Range("D3").Value = "'1.810.24"
$ws.Range("D5").Value = "'226.03"
$ws.Range("D6").Value = "'0.600"
$ws.Range("D8").Value = "'36.34"
$ws.Range("D11").Value = "'0.0967"
$ws.Range("D13").Value = "'11.31"
$ws.Range("D14").Value = "'1.843.92"
$ws.Range("D15").Value = "'0.631"
$ws.Range("D16").Value = "'34.461.62"
$ws.Range("D17").Value = "'4.43"
$ws.Range("D18").Value = "'68.71"
$ws.Range("D25").Value = "'171.51"
$ws.Range("D26").Value = "'7.95"
$ws.Range("D27").Value = "'17.27"
$ws.Range("D28").Value = "'0.122"
$ws.Range("D30").Value = "'3.83"
$ws.Range("D31").Value = "'3.94"
$ws.Range("D33").Value = "'0.0519"
$ws.Range("D34").Value = "'1.81"
$ws.Range("D35").Value = "'1.363.62"
$ws.Range("D36").Value = "'0.654"
$ws.Range("D39").Value = "'0.0187"
$ws.Range("D40").Value = "'2.43"
$ws.Range("D41").Value = "'2.79"
$ws.Range("D42").Value = "'81.12"
$ws.Range("D43").Value = "'0.938"
$ws.Range("D44").Value = "'1.16"
$ws.Range("D45").Value = "'13.39"
$ws.Range("D46").Value = "'0.0498"
$ws.Range("D47").Value = "'1.971.29"
$ws.Range("D50").Value = "'102.81"
$ws.Range("D51").Value = "'0.0" + [char]0x2086 + "0123"
$ws.Range("D2,D3,D5,D6,D8,D11,D13,D14,D15,D16,D17,D18,D25,D26,D27,D28,D30,D31,D33,D34,D35,D36,D39,D40,D41,D42,D43,D44,D45,D46,D47,D50,D51").Style = "Normal"

# --- 1h volume-change column (E): text already contains "%" and padding
#     spaces, so a plain assignment keeps it as a string as-is.
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  -2.23%  "
$ws.Range("E10").Value = "  -1.88%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -1.27%  "
$ws.Range("E24").Value = "  +5.09%  "
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("E27").Value = "  +2.93%  "
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -1.68%  "
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  -2.47%  "
$ws.Range("E36").Value = "  -4.03%  "
$ws.Range("E37").Value = "  -0.58%  "
$ws.Range("E38").Value = "  -5.61%  "
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("E44").Value = "  +4.26%  "
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("E48").Value = "  -2.57%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("E51").Value = "  -5.23%  "
